$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExcelGuru99Demo")

# Row 1 - header-like row
$ws.Range("A1").Value = "mukesh"
$ws.Range("B1").Value = "kjjhjkhkj"
$ws.Range("C1").Value = "rrtfdf"
$ws.Range("D1").Value = "rt4tew"
$ws.Range("E1").Value = "tretertr"
$ws.Range("F1").Value = "retretret"
$ws.Range("G1").Value = "retret"
$ws.Range("H1").Value = "retertert"

# Row 2
$ws.Range("A2").Value = "priyanks"
$ws.Range("B2").Value = "dfdfdsf"
$ws.Range("C2").Value = "fddsfdf"

# Row 3
$ws.Range("A3").Value = "weewrew"
$ws.Range("B3").Value = "dfdsfdsf"
$ws.Range("C3").Value = "fdsfsdf"

# Row 4
$ws.Range("A4").Value = "retert"
$ws.Range("B4").Value = "dfdsfdsf"
$ws.Range("C4").Value = "dfdsf"

# Row 5
$ws.Range("A5").Value = "reegret"
$ws.Range("B5").Value = "dfdsfdsf"
$ws.Range("C5").Value = "dsfsf"

$ws.Range("G1").Select() | Out-Null
